# Update countries & provincias Spain
# - Swap Chile/Singapur ranking rows (30/31) with refreshed case counts
# - Swap Estonia/Armenia ranking rows (71/72) with refreshed case counts
# - Bump the "Datos actualizados" timestamp from 09:22 to 09:52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer note (row 1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 09:52"

# Rows 30/31: Chile and Singapur swap places in the ranking
$ws.Range("A30").Value = "Singapur"
$ws.Range("B30").Value = 12693
$ws.Range("C30").Value = 618
$ws.Range("D30").Value = 956
$ws.Range("E30").Value = 11725
$ws.Range("F30").Value = 24
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 12

$ws.Range("A31").Value = "Chile"
$ws.Range("B31").Value = 12306
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 6327
$ws.Range("E31").Value = 5805
$ws.Range("F31").Value = 408
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 174

# Rows 71/72: Estonia and Armenia swap places in the ranking
$ws.Range("A71").Value = "Armenia"
$ws.Range("B71").Value = 1677
$ws.Range("C71").Value = 81
$ws.Range("D71").Value = 803
$ws.Range("E71").Value = 846
$ws.Range("F71").Value = 10
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 28

$ws.Range("A72").Value = "Estonia"
$ws.Range("B72").Value = 1605
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 206
$ws.Range("E72").Value = 1353
$ws.Range("F72").Value = 6
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 46
